$d = $word.ActiveDocument

function Get-ParaIndexContainingText {
    param([string]$searchText)
    $r = $d.Range(0, 0)
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return -1 }
    $pos = $r.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# Work from the bottom of the document upward so that earlier
# (text-anchor based) lookups remain valid while later parts of the
# document are edited first.
# ------------------------------------------------------------------

# 5) After the "IF we have time..." paragraph, add a new bold
#    paragraph "Ran out of time".
$idxIfWeHaveTime = Get-ParaIndexContainingText("IF we have time, go over the schedule")
$pIfWeHaveTime = $d.Paragraphs.Item($idxIfWeHaveTime)
$pIfWeHaveTime.Range.InsertAfter("`rRan out of time")

# 4) After the "Schedule an evening meeting..." paragraph: delete the
#    blank paragraph that follows it, then add a new bold paragraph
#    "TONIGHT 4pm" followed by a blank paragraph, before "IF we have
#    time...".
$idxSchedule = Get-ParaIndexContainingText("Schedule an evening meeting")
$pBlankAfterSchedule = $d.Paragraphs.Item($idxSchedule + 1)
$pBlankAfterSchedule.Range.Delete()

$pSchedule = $d.Paragraphs.Item($idxSchedule)
$pSchedule.Range.InsertAfter("`rTONIGHT 4pm`r")

# 3) The blank paragraph right after "Wrap up of the presentation..."
#    becomes a new note paragraph (plain, not bold).
$idxWrapUp = Get-ParaIndexContainingText("Wrap up of the presentation")
$pBlankAfterWrapUp = $d.Paragraphs.Item($idxWrapUp + 1)
$pBlankAfterWrapUp.Range.InsertAfter("Note taken, people need to split up sections (some are too long, some are too short) and condense info on the slides specifically the requirement slides and probably some of the low level stuff.")

# 2) Insert the meeting-minutes block right after the first blank
#    paragraph that follows "Mock presentation 40 mins" (i.e. before
#    the second blank paragraph that precedes "Wrap up...").
$idxMock = Get-ParaIndexContainingText("Mock presentation 40 mins")
$pBlankAfterMock = $d.Paragraphs.Item($idxMock + 1)

$minutesBlock = @(
    "Connor:",
    "5 mins ended at batched order from the requirements from Connors part",
    "Requirements ended up being 14mins",
    "Merchant and supplier both have order tracking, we can condense this down into 1 slide and mention that is applies to the two parties",
    "",
    "Josh:",
    "Repeated what I said at the beginning with the users",
    "3 mins",
    "",
    "Lincoln:",
    "2.5 mins",
    "",
    "Galmo:",
    "Inconsistencies with lincolns (in transit) wording",
    "10 mins",
    "",
    "Micheal:",
    "Went short because we were short on time",
    "His part went like 2 mins",
    "",
    "Shiva:",
    "1.5 mins",
    "",
    "",
    "",
    ""
) -join "`r"

$pBlankAfterMock.Range.InsertAfter("`r" + $minutesBlock)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
